$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(8106.63, 8034.32, 18.84, 19.010000000000002, $false, 0.9,     42613.765567129631, $true),
    @(8062.04, 8106.63, 18.93, 18.824999999999999, $false, -0.55000000000000004, 42614.67292824074, $false),
    @(8019.31, 8062.04, 18.72, 18.62,               $false, -0.53,  42615.750173611108, $false)
)

$r = 10
foreach ($row in $data) {
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
